$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hoja.tareas")

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "gghg"
$ws.Range("C2").Value = "tfyty"
$ws.Range("E2").Value = "ftytg"
$ws.Range("F2").Value = "tgyty"
$ws.Range("D2").Value = "pendi"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "tff"
$ws.Range("C3").Value = "gfgdf"
$ws.Range("D3").Value = "eje"
$ws.Range("E3").Value = "gfgd"
$ws.Range("F3").Value = "dfg"

$ws.Range("H8").Select()
